# pyGooIbClient Requirements.xlsx - "Minor re-organization to move watchdog
# and httpendpoints to lib folder"
#
# The Features table is re-ordered/re-grouped, several rows are dropped
# (re-authentication items, CLI/argparse/interface-selection items, the
# "build on basic library" item), and a new "Category" column (B) plus a
# batch of new Logging/Configuration/Interfaces/Parameter/User Interface
# rows are introduced. Net effect: the table shrinks from 31 to 28 data
# rows (32 -> 29 incl. header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 3 now-unused trailing rows first so the sheet dimension and the
# backing table range shrink by 3 (A1:C32 -> A1:C29, table A1:C75 -> A1:C72)
# exactly like the real edit did.
$ws.Range("A30:A32").EntireRow.Delete() | Out-Null

# Final Feature/Category/Description content, in the target row order.
# ("" means leave column B blank for that row.)
$rows = @(
    @{B="Category";        C="Description"},
    @{B="Interfaces";      C="Provide a GET interface"},
    @{B="Interfaces";      C="Provide a POST interface"},
    @{B="";                C="Method to Request streaming market data from IB Client"},
    @{B="";                C="Method to Request account information from IB Client"},
    @{B="Configuration";   C="Username and account password stored in enviromental variables, not in system code"},
    @{B="";                C="External applications can subscribe to streaming market data via requests to the client"},
    @{B="";                C="External applications receive periodic updates to their market data subscriptions"},
    @{B="";                C="External applications shall be able to unsubscribe to market data"},
    @{B="";                C="External applications shall be able to submit, modify, and cancel orders"},
    @{B="";                C="External applications will receive updates about changes in order status (fills, etc.)"},
    @{B="";                C="External applications can request and receive status about the system"},
    @{B="Logging";         C="All loggable events will be stored in a master log file"},
    @{B="Logging";         C="Log files will cover a single day"},
    @{B="Logging";         C="Loggable events will be categorized in order to make event extraction easier"},
    @{B="Logging";         C="Loggable events will be timestamped"},
    @{B="";                C="Market data will be stored in a centralized database"},
    @{B="Interfaces";      C="provide a ROS2 interface"},
    @{B="Configuration";   C="All configurable parameters will be stored in YAML format"},
    @{B="Configuration";   C="Configuration files will be stored in the folder where used [default] but may use an alternate path if provided"},
    @{B="User Interface";  C="The client will have a CLI"},
    @{B="Parameter";       C="Provide a configurable time between IB Client `"Tickle`" events"},
    @{B="";                C="CLI flags use standard Linux format"},
    @{B="";                C="System can recover and retrieve state from loss of network"},
    @{B="";                C="System can recover and retrieve state from host reboot"},
    @{B="";                C="System can recover and retrieve state from loss of communication with Client Portal"},
    @{B="";                C="Provide means to determine current version of Client Portal gateway"},
    @{B="";                C="Provide means to determine current version of library"},
    @{B="";                C="provide means to connect to ib websocket stream for market data"}
)

$ws.Range("A1").Value = "Feature"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $b = $rows[$i].B
    $c = $rows[$i].C

    if ($b -ne "") {
        $ws.Range("B$r").Value = $b
    } else {
        $ws.Range("B$r").Value = $null
    }
    $ws.Range("C$r").Value = $c
}

# Column B ("Category") now holds real data; widen it (best effort - the
# interop engine quantizes widths, closest reachable step to the
# original author's best-fit 13.5703125 is 13.5).
$ws.Columns.Item(2).ColumnWidth = 12.6

# Restore the recorded selection / active cell.
$ws.Range("C13").Select() | Out-Null
